$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.161928772926331
$ws.Range("B1").Value = 2.414951801300049
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.381839752197266
$ws.Range("E1").Value = 1.233621716499329
